$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q2").Value = 625526
$ws.Range("R2").Value = 6542585

$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
